$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray volatile-formula row (C50) that is no longer needed.
$ws.Rows("50:50").Delete()

# Add the new "Car Fleet" problem entry as row 44, matching the look of
# the row above it (row 43).
$ws.Range("A43").Copy()
$ws.Range("A44").PasteSpecial(-4122)

$ws.Range("A44").Value2 = "Car Fleet"
$ws.Range("D44").Value2 = "https://leetcode.com/problems/car-fleet/"

# Re-create the hyperlink for the new URL cell, then restore the proper
# "Hyperlink" cell style (Hyperlinks.Add re-styles the cell on its own).
$ws.Hyperlinks.Add($ws.Range("D44"), "https://leetcode.com/problems/car-fleet/")
$ws.Range("D43").Copy()
$ws.Range("D44").PasteSpecial(-4122)

# Update the selection / scroll position shown when the sheet is opened.
$ws.Range("C27").Select()
